# Add four new "Title and Content" slides (Week 1 - Tuesday material, plus the
# updated Monday Chapter 3 slide) at the end of the deck, mirroring the layout
# already used by the existing slide 4 (sldId 259).

$p = $ppt.ActivePresentation

# "Title and Content" is the 2nd custom layout on the slide master (same
# layout slide 3 & slide 4 already use).
$layout = $p.SlideMaster.CustomLayouts.Item(2)

$vt = [char]11   # vertical-tab -> soft line break between title lines
$nl = [char]13   # carriage return -> new bullet paragraph in body text

function Set-BodyText($shape, [string[]]$paragraphs) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = [string]::Join($nl, $paragraphs)
    $tr.LanguageID = "nl-NL"
}

function Set-TitleText($shape, [string]$line1, [string]$line2) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = $line1 + $vt + $line2
    $tr.LanguageID = "nl-NL"
}

# ---------------------------------------------------------------------------
# Slide 5 (sldId 260) - Week 1 - Monday / Chapter 3 (updated exercise slide)
# ---------------------------------------------------------------------------
$s5 = $p.Slides.AddSlide(5, $layout)

Set-TitleText $s5.Shapes.Item(1) "Week 1 – Monday" "Chapter 3:"

$s5body = $s5.Shapes.Item(2)
Set-BodyText $s5body @(
    "Inheriting a class variable shares it with all instances of the class, and there is only one instance of the ParentClass.classvariable.",
    "You  can extend any built-in class with inheritence, like list, set, dict, file, str, int, float.",
    "Methods can be overwritten, even __init__",
    "Use super() to get an instance of the parent object.",
    "Avoid multiple inheritance, it’s more trouble than it’s worth.",
    "Polymorphism -> extend parent class with inheritence. Parent class can check functions from subclasses without knowing which subclass it is referring to.",
    "Duck typing: a type or class of an object is less important than the method it defines -> can create relationships without setting up code for inheritence.",
    "Abstract base class: define set of methods and properties that a base class must implement.",
    "",
    ""
)
$s5body.TextFrame.AutoSize = 3

# ---------------------------------------------------------------------------
# Slide 6 (sldId 261) - Week 1 - Tuesday / Chapter 4
# ---------------------------------------------------------------------------
$s6 = $p.Slides.AddSlide(6, $layout)

Set-TitleText $s6.Shapes.Item(1) "Week 1 – Tuesday" "Chapter 4:"

$s6body = $s6.Shapes.Item(2)
Set-BodyText $s6body @(
    "Exceptions raised will stop the function from continuing onwards unless handled (try, except).",
    "Using ‘except EceptionName:’ catches specific exceptions. Don’t use ‘except’ on its own, since it will als catch SystemExit and KeyboardInterrupt.",
    "Can use ‘except(Error1, Error2, ...)’ to catch multiple exceptions with one statment.",
    "Stacking exception clauses will only cause the first one to execute.",
    "‘raise’ re-raises last exception.",
    "Can use ‘as’ to capture exception as variable and use it after.",
    "Use ‘else’ after ‘except’ in case no exception is called and use ‘finally’ after ‘except’ to run code no matter what exception is called (think: closing an open file).",
    "‘finally’ will also occur before a return in a try clause."
)
$s6body.TextFrame.AutoSize = 3

# ---------------------------------------------------------------------------
# Slide 7 (sldId 262) - Week 1 - Tuesday / Chapter 4 (continued)
# ---------------------------------------------------------------------------
$s7 = $p.Slides.AddSlide(7, $layout)

Set-TitleText $s7.Shapes.Item(1) "Week 1 – Tuesday" "Chapter 4:"

$s7body = $s7.Shapes.Item(2)
Set-BodyText $s7body @(
    "SystemExit and KeyboardInterrupt inherit from BaseException instead of Exception so we can clean up before the program exits.",
    "Can create new errors as a class, then use ‘except errorname as e’ to use ‘e.errorclassmethod’ in handler. (For example, to return a bank balance).",
    "Python programmers tend to execute code and deal with what goes wrong with exceptions, instead of planning for every thing that might happen, but is not going to arise in the normal path through code."
)
$s7body.TextFrame.AutoSize = 3

# ---------------------------------------------------------------------------
# Slide 8 (sldId 263) - Week 1 - Tuesday / PEP8
# ---------------------------------------------------------------------------
$s8 = $p.Slides.AddSlide(8, $layout)

Set-TitleText $s8.Shapes.Item(1) "Week 1 – Tuesday" "PEP8:"

$s8body = $s8.Shapes.Item(2)
Set-BodyText $s8body @(
    "Know when to ignore the guideline (makes code less readable, breaks backwards compatibility, etc.).",
    "Indent = 4 spaces.",
    "Max line length = 79. 72 for comments and docstrings.",
    "Top level functions + classes 2 blank lines, method 1 blank line.",
    "UTF-8",
    "Only use inline comments if they actually add value.",
    "Use Pylint to easily follow PEP8 coding standards."
)
$s8body.TextFrame.AutoSize = 3

Write-Host "Added slides: $($p.Slides.Count) total"
